$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# Cells holding plain "NN%" text need an explicit Text format first,
# otherwise Excel auto-converts the string into a numeric percentage.
$ws.Range('H12').NumberFormat = '@'
$ws.Range('H13').NumberFormat = '@'
$ws.Range('H16').NumberFormat = '@'
$ws.Range('H21').NumberFormat = '@'
$ws.Range('H25').NumberFormat = '@'
$ws.Range('H29').NumberFormat = '@'
$ws.Range('H30').NumberFormat = '@'
$ws.Range('H36').NumberFormat = '@'
$ws.Range('H38').NumberFormat = '@'
$ws.Range('H40').NumberFormat = '@'

$ws.Range('E2').Value = '2026-03-01 06:18:40'
$ws.Range('N2').Value = '-2.5 °C 5:47 TU'
$ws.Range('O2').Value = '-1.0 °C'
$ws.Range('E3').Value = '2026-03-01 06:18:43'
$ws.Range('L3').Value = '20.5 km/h - 123º 5:30 TU'
$ws.Range('N3').Value = '-4.2 °C 5:46 TU'
$ws.Range('O3').Value = '-3.7 °C'
$ws.Range('E4').Value = '2026-03-01 06:18:46'
$ws.Range('E5').Value = '2026-03-01 06:18:48'
$ws.Range('N5').Value = '-5.0 °C 5:44 TU'
$ws.Range('O5').Value = '-3.8 °C'
$ws.Range('E6').Value = '2026-03-01 06:18:51'
$ws.Range('E7').Value = '2026-03-01 06:18:54'
$ws.Range('E8').Value = '2026-03-01 06:18:56'
$ws.Range('J8').Value = '1025.6 hPa'
$ws.Range('E9').Value = '2026-03-01 06:18:59'
$ws.Range('O9').Value = '11.6 °C'
$ws.Range('E10').Value = '2026-03-01 06:19:02'
$ws.Range('K10').Value = '-0.1 MJ/m2'
$ws.Range('N10').Value = '3.9 °C 5:59 TU'
$ws.Range('O10').Value = '6.1 °C'
$ws.Range('E11').Value = '2026-03-01 06:19:04'
$ws.Range('N11').Value = '5.8 °C 5:59 TU'
$ws.Range('E12').Value = '2026-03-01 06:19:07'
$ws.Range('H12').Value = '74%'
$ws.Range('N12').Value = '8.5 °C 5:30 TU'
$ws.Range('O12').Value = '10.3 °C'
$ws.Range('E13').Value = '2026-03-01 06:19:09'
$ws.Range('H13').Value = '92%'
$ws.Range('N13').Value = '3.2 °C 5:57 TU'
$ws.Range('E14').Value = '2026-03-01 06:19:12'
$ws.Range('N14').Value = '9.2 °C 5:42 TU'
$ws.Range('O14').Value = '10.8 °C'
$ws.Range('E15').Value = '2026-03-01 06:19:14'
$ws.Range('E16').Value = '2026-03-01 06:19:17'
$ws.Range('H16').Value = '84%'
$ws.Range('N16').Value = '-6.8 °C 5:59 TU'
$ws.Range('O16').Value = '-5.1 °C'
$ws.Range('E17').Value = '2026-03-01 06:19:20'
$ws.Range('E18').Value = '2026-03-01 06:19:22'
$ws.Range('J18').Value = '1025.9 hPa'
$ws.Range('E19').Value = '2026-03-01 06:19:24'
$ws.Range('I19').Value = '0.2 mm'
$ws.Range('N19').Value = '5.9 °C 5:45 TU'
$ws.Range('E20').Value = '2026-03-01 06:19:27'
$ws.Range('L20').Value = '10.8 km/h - 202º 5:58 TU'
$ws.Range('N20').Value = '-4.6 °C 5:57 TU'
$ws.Range('O20').Value = '-3.1 °C'
$ws.Range('E21').Value = '2026-03-01 06:19:29'
$ws.Range('H21').Value = '88%'
$ws.Range('J21').Value = '1025.5 hPa'
$ws.Range('L21').Value = '5.8 km/h - 234º 5:49 TU'
$ws.Range('N21').Value = '5.2 °C 5:59 TU'
$ws.Range('O21').Value = '6.4 °C'
$ws.Range('E22').Value = '2026-03-01 06:19:32'
$ws.Range('L22').Value = '14.8 km/h - 325º 5:35 TU'
$ws.Range('N22').Value = '-6.5 °C 5:53 TU'
$ws.Range('O22').Value = '-5.4 °C'
$ws.Range('E23').Value = '2026-03-01 06:19:35'
$ws.Range('O23').Value = '-3.6 °C'
$ws.Range('E24').Value = '2026-03-01 06:19:38'
$ws.Range('O24').Value = '4.6 °C'
$ws.Range('E25').Value = '2026-03-01 06:19:40'
$ws.Range('H25').Value = '94%'
$ws.Range('E26').Value = '2026-03-01 06:19:43'
$ws.Range('N26').Value = '2.3 °C 5:31 TU'
$ws.Range('E27').Value = '2026-03-01 06:19:45'
$ws.Range('E28').Value = '2026-03-01 06:19:48'
$ws.Range('E29').Value = '2026-03-01 06:19:50'
$ws.Range('H29').Value = '80%'
$ws.Range('K29').Value = '-0.1 MJ/m2'
$ws.Range('N29').Value = '6.7 °C 5:54 TU'
$ws.Range('O29').Value = '9.1 °C'
$ws.Range('E30').Value = '2026-03-01 06:19:53'
$ws.Range('H30').Value = '79%'
$ws.Range('E31').Value = '2026-03-01 06:19:56'
$ws.Range('J31').Value = '1024.6 hPa'
$ws.Range('E32').Value = '2026-03-01 06:19:58'
$ws.Range('M32').Value = '5.0 °C 5:51 TU'
$ws.Range('O32').Value = '2.7 °C'
$ws.Range('E33').Value = '2026-03-01 06:20:01'
$ws.Range('N33').Value = '3.6 °C 5:59 TU'
$ws.Range('E34').Value = '2026-03-01 06:20:04'
$ws.Range('N34').Value = '-0.8 °C 5:52 TU'
$ws.Range('E35').Value = '2026-03-01 06:20:06'
$ws.Range('E36').Value = '2026-03-01 06:20:09'
$ws.Range('H36').Value = '75%'
$ws.Range('L36').Value = '21.6 km/h - 354º 5:59 TU'
$ws.Range('O36').Value = '10.1 °C'
$ws.Range('E37').Value = '2026-03-01 06:20:12'
$ws.Range('N37').Value = '6.1 °C 5:42 TU'
$ws.Range('E38').Value = '2026-03-01 06:20:14'
$ws.Range('H38').Value = '97%'
$ws.Range('O38').Value = '8.7 °C'
$ws.Range('E39').Value = '2026-03-01 06:20:17'
$ws.Range('E40').Value = '2026-03-01 06:20:19'
$ws.Range('H40').Value = '89%'
$ws.Range('N40').Value = '5.4 °C 5:44 TU'
$ws.Range('O40').Value = '6.9 °C'
$ws.Range('E41').Value = '2026-03-01 06:20:21'
$ws.Range('J41').Value = '1025.5 hPa'
$ws.Range('N41').Value = '10.8 °C 5:59 TU'
$ws.Range('E42').Value = '2026-03-01 06:20:23'
$ws.Range('N42').Value = '5.2 °C 5:56 TU'
$ws.Range('O42').Value = '8.3 °C'
$ws.Range('E43').Value = '2026-03-01 06:20:26'
$ws.Range('E44').Value = '2026-03-01 06:20:29'
$ws.Range('N44').Value = '-4.3 °C 5:48 TU'
$ws.Range('O44').Value = '-2.8 °C'
$ws.Range('E45').Value = '2026-03-01 06:20:31'
$ws.Range('G45').Value = '2 cm'
$ws.Range('N45').Value = '2.0 °C 5:51 TU'
$ws.Range('O45').Value = '3.4 °C'
$ws.Range('E46').Value = '2026-03-01 06:20:34'
